$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells whose new text values look like plain numbers need the Text
# number format forced first so Excel keeps them as literal strings
# (matching the original inline-string cells) instead of converting
# them to numeric values and losing formatting such as trailing zeros.

$ws.Range('D2').Value = '24.329.77'
$ws.Range('E2').Value = '  +1.00%  '

$ws.Range('D3').Value = '1.668.15'
$ws.Range('E3').Value = '  +1.66%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.89'
$ws.Range('E5').Value = '  +1.19%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.13%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3953'
$ws.Range('E7').Value = '  +1.69%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3932'
$ws.Range('E8').Value = '  +1.77%  '

$ws.Range('E9').Value = '  +5.10%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.389'
$ws.Range('E10').Value = '  +3.09%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.002'
$ws.Range('E11').Value = '  -0.09%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08577'
$ws.Range('E12').Value = '  -0.95%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.41'
$ws.Range('E13').Value = '  +3.53%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.311'
$ws.Range('E14').Value = '  +3.09%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.037'
$ws.Range('E15').Value = '  +7.81%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001334'
$ws.Range('E16').Value = '  +3.71%  '

$ws.Range('D17').Value = '1.670.46'
$ws.Range('E17').Value = '  +2.04%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '95.51'
$ws.Range('E18').Value = '  +0.18%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07028'
$ws.Range('E19').Value = '  +1.69%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.48'
$ws.Range('E20').Value = '  -0.23%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.000'
$ws.Range('E21').Value = '  +1.54%  '

$ws.Range('E22').Value = '  +0.22%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.78'
$ws.Range('E23').Value = '  +1.69%  '

$ws.Range('D24').Value = '24.341.51'
$ws.Range('E24').Value = '  +1.09%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.563'
$ws.Range('E25').Value = '  +9.93%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.118'
$ws.Range('E26').Value = '  +13.21%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.50'
$ws.Range('E27').Value = '  +0.67%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.21'
$ws.Range('E28').Value = '  -0.50%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '142.65'
$ws.Range('E29').Value = '  +1.54%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.371'
$ws.Range('E30').Value = '  +1.27%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.036'
$ws.Range('E31').Value = '  -5.24%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.536'
$ws.Range('E32').Value = '  +4.72%  '

$ws.Range('D33').Value = '1.848.63'
$ws.Range('E33').Value = '  +2.15%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.063'
$ws.Range('E34').Value = '  +12.51%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.03106'
$ws.Range('E35').Value = '  +7.40%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08234'
$ws.Range('E36').Value = '  +2.78%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.890'
$ws.Range('E37').Value = '  -0.14%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '11.21'
$ws.Range('E38').Value = '  +13.43%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2762'
$ws.Range('E39').Value = '  +3.25%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09287'
$ws.Range('E40').Value = '  +0.94%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.7692'
$ws.Range('E41').Value = '  +2.13%  '

$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.77'
$ws.Range('E42').Value = '  +5.97%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.440'
$ws.Range('E43').Value = '  -1.86%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.55'
$ws.Range('E44').Value = '  +3.78%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7081'
$ws.Range('E45').Value = '  +2.88%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.536'
$ws.Range('E46').Value = '  +3.14%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.121'
$ws.Range('E47').Value = '  +0.93%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').Value = '  +0.12%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08423'
$ws.Range('E49').Value = '  +0.34%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '136.65'
$ws.Range('E50').Value = '  +2.86%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.264'
$ws.Range('E51').Value = '  +0.61%  '
